$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated sval data (filtering save games) - updated B:E and G (sum) columns for rows 2-8.
# Column F (Win) and column A (date) are unchanged.

$data = @{
    2 = @{ B = 3.182878228561681;  C = 1.65323645889881;  D = 0.1529057820181812;  E = 0.4998867070740569; G = 5.488907176552729 }
    3 = @{ B = 3.182878228561681;  C = 1.65323645889881;  D = 0.7127328510149897;  E = 6.48142807727062;   G = 12.0302756157461 }
    4 = @{ B = 3.182878228561681;  C = 1.65323645889881;  D = 3.082599426703578;   E = 6.48142807727062;   G = 14.40014219143469 }
    5 = @{ B = 0.7287194209349384; C = 1.65323645889881;  D = 3.082599426703578;   E = 0.4998867070740569; G = 5.964442013611383 }
    6 = @{ B = 1.505614041169197;  C = 1.65323645889881;  D = 0.1529057820181812;  E = 0.4998867070740569; G = 3.811642989160245 }
    7 = @{ B = 0.7287194209349384; C = 1.65323645889881;  D = 0.1529057820181812;  E = 0.4998867070740569; G = 3.034748368925986 }
    8 = @{ B = 3.182878228561681;  C = 1.65323645889881;  D = 16.98373111632243;   E = 6.48142807727062;   G = 28.30127388105354 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
